# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table (rows 3-5) ---
$ws.Range("C3").Value = 493
$ws.Range("D3").Value = 90.2
$ws.Range("D4").Value = 98.59999999999999
$ws.Range("C5").Value = 555

# --- "Good Drivers" table (rows 13-18) ---
# The six driver rows get re-sorted (by sample count, descending) and the
# "Driver Vintage" dates get refreshed for this week's run.

# New row 13 now holds the data that used to be driver 23.100.0.4
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B13").Value = 445055
$ws.Range("D13").Value = 99.90000000000001

# New row 14 now holds the data that used to be driver 22.80.0.9
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B14").Value = 77849
$ws.Range("D14").Value = 99.90000000000001

# New row 15 now holds the data that used to be driver 22.50.1.1
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B15").Value = 34244
$ws.Range("D15").Value = 100

# New row 16 now holds the data that used to be driver 21.110.3.2
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B16").Value = 59673
$ws.Range("D16").Value = 100

# New row 17 now holds the data that used to be driver 21.70.0.6
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B17").Value = 113652
$ws.Range("D17").Value = 100

# New row 18 now holds the data that used to be driver 21.60.2.1
# (its "Driver Vintage" text, E18, is untouched this week)
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B18").Value = 56018
$ws.Range("D18").Value = 100

# "Driver Vintage" (column E) text values. Pre-format as Text so Excel does
# not auto-convert the yyyy-mm-dd-looking strings into date serials, then
# copy the original "General" number format back on top so the stored
# style index matches the rest of the column.
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2024-11-10"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2021-08-18"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2021-04-27"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2020-08-05"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2020-01-06"

# E18 ("2019-12-14") is unchanged this week.

$ws.Range("D17").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E17").PasteSpecial(-4122)
